# Update "想去人数" (want-to-go count, column F) figures for a handful of
# events, on both the "展览" sheet and the "全部类型" sheet (which mirrors
# the same events plus one extra row, so the rows are offset by one there).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3272
$ws1.Range("F5").Value = 6905
$ws1.Range("F6").Value = 2176
$ws1.Range("F7").Value = 34
$ws1.Range("F13").Value = 157

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3272
$ws4.Range("F6").Value = 6905
$ws4.Range("F7").Value = 2176
$ws4.Range("F8").Value = 34
$ws4.Range("F14").Value = 157
